$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.205588221549988
$ws.Range("B1").Value = 2.323080062866211
$ws.Range("C1").Value = 4.672676563262939
$ws.Range("D1").Value = 3.376785755157471
$ws.Range("E1").Value = 1.177467107772827
